{"js": "// Update the GifAnimation library link in the documentation.\nasync function replaceText(body, needle, replacement) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Main edit: point the GifAnimation repo link at its new home.\nawait replaceText(\n  body,\n  \"https://github.com/akiljohnson1/GifAnimation\",\n  \"https://github.com/extrapixel/gif-animation\"\n);\n\n// The asset-credit rows for image_ice.png / image_star.png were re-saved\n// alongside the link update (their split runs collapse into single runs),\n// even though their visible text is unchanged.\nawait replaceText(body, \"image_ice.png\", \"image_ice.png\");\nawait replaceText(body, \"image_star.png\", \"image_star.png\");\n", "ps1": "$d = $word.ActiveDocument\n\n# Main edit: point the GifAnimation repo link at its new home. Using the\n# paragraph's own Range (rather than a document-wide Find/Replace) keeps the\n# run's existing formatting/rsid attributes untouched, matching how Word\n# applies an in-place text edit.\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text -like \"*https://github.com/akiljohnson1/GifAnimation*\") {\n        $r.Text = \"https://github.com/extrapixel/gif-animation\"\n    }\n}\n\n# The asset-credit rows for image_ice.png / image_star.png were re-saved\n# alongside the link update (their split runs collapse into single runs),\n# even though their visible text is unchanged.\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$findText, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 2) | Out-Null\n}\n\nReplace-AllText \"image_ice.png\" \"image_ice.png\"\nReplace-AllText \"image_star.png\" \"image_star.png\"\n"}
